# urto anaelastico con pesi
# Applies:
#  1) "urti-elastici": fix sign of column A (now positive), correct two B values,
#     and add Eki/Ekf/deltaEk + ViB/Vf/deltaP derived columns (E..J).
#  2) "urti-anaelastici": fill in A/B measured data and add the full set of
#     derived columns (F..M): Eki/Ekf/deltaEk, deltaV, deltaV/Vi, Pi/Pf/deltaP.
#  3) Two brand-new worksheets "urto-anaelastico+massa" and
#     "urto-elastico+massa" with their starter data/headers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) urti-elastici
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("urti-elastici")

# Column A sign flip (measured v1 was recorded with the wrong sign).
$ws.Range("A2").Value = 0.449
$ws.Range("A3").Value = 0.63
$ws.Range("A4").Value = 0.624
$ws.Range("A5").Value = 0.57
$ws.Range("A6").Value = 0.588
$ws.Range("A7").Value = 0.809
$ws.Range("A8").Value = 0.863
$ws.Range("A9").Value = 0.663

# Corrected B (vf) readings for rows 7 & 8.
$ws.Range("B7").Value = 0.837
$ws.Range("B8").Value = 0.878

# New header row (Eki, Ekf, deltaEk, Pi, Pf, deltaP).
$ws.Range("E1").Value = "Eki"
$ws.Range("F1").Value = "Ekf"
$ws.Range("G1").Value = "deltaEk"
$ws.Range("H1").Value = "Pi"
$ws.Range("I1").Value = "Pf"
$ws.Range("J1").Value = "deltaP"

for ($r = 2; $r -le 9; $r++) {
    $ws.Range("E$r").Formula = "=0.5*(`$D`$2)*B$r*B$r"
    $ws.Range("F$r").Formula = "=0.5*`$C`$2*A$r*A$r"
    $ws.Range("G$r").Formula = "=F$r-E$r"
    $ws.Range("H$r").Formula = "=`$D`$2*B$r"
    $ws.Range("I$r").Formula = "=`$C`$2*A$r"
    $ws.Range("J$r").Formula = "=I$r-H$r"
    $ws.Range("E$r`:F$r").NumberFormat = "0.00000"
    $ws.Range("H$r`:I$r").NumberFormat = "0.00000"
    $ws.Range("G$r").NumberFormat = "0.000"
    $ws.Range("J$r").NumberFormat = "0.000"
}

# ---------------------------------------------------------------------------
# 2) urti-anaelastici
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("urti-anaelastici")

$aVals = @(-0.37, -0.502, -0.615, -0.656, -0.547, -0.679, -0.556, -0.535, -0.556)
$bVals = @(-0.178, -0.243, -0.302, -0.32, -0.264, -0.33, -0.269, -0.26, -0.27)
for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2
    $ws2.Range("A$r").Value = $aVals[$i]
    $ws2.Range("B$r").Value = $bVals[$i]
}

# New header row (ViB, Vf, mR, mB, ViR, Eki, Ekf, deltaEk, deltaV, deltaV/Vi, Pi, Pf, deltaP).
$ws2.Range("A1").Value = "ViB"
$ws2.Range("B1").Value = "Vf"
$ws2.Range("C1").Value = "mR"
$ws2.Range("D1").Value = "mB"
$ws2.Range("E1").Value = "ViR"
$ws2.Range("F1").Value = "Eki"
$ws2.Range("G1").Value = "Ekf"
$ws2.Range("H1").Value = "deltaEk"
$ws2.Range("I1").Value = "deltaV"
$ws2.Range("J1").Value = "deltaV/Vi"
$ws2.Range("K1").Value = "Pi"
$ws2.Range("L1").Value = "Pf"
$ws2.Range("M1").Value = "deltaP"

for ($r = 2; $r -le 10; $r++) {
    $ws2.Range("F$r").Formula = "=0.5*`$D`$2*A$r*A$r"
    $ws2.Range("G$r").Formula = "=0.5*(`$C`$2+`$D`$2)*B$r*B$r"
    $ws2.Range("H$r").Formula = "=G$r-F$r"
    $ws2.Range("I$r").Formula = "=B$r-A$r"
    $ws2.Range("J$r").Formula = "=I$r/A$r"
    $ws2.Range("K$r").Formula = "=`$D`$2*A$r"
    $ws2.Range("L$r").Formula = "=(`$C`$2+`$D`$2)*B$r"
    $ws2.Range("M$r").Formula = "=L$r-K$r"

    $ws2.Range("F$r`:G$r").NumberFormat = "0.00000"
    $ws2.Range("H$r").NumberFormat = "0.000"
    $ws2.Range("I$r").NumberFormat = "0.000"
    $ws2.Range("K$r`:L$r").NumberFormat = "0.000"
    $ws2.Range("M$r").NumberFormat = "0.000"

    $ws2.Range("J$r").NumberFormat = "0.00"
    $ws2.Range("J$r").Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 3) New worksheet "urto-anaelastico+massa"
# ---------------------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIdx)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "urto-anaelastico+massa"

$ws3.Range("A1").Value = "viB"
$ws3.Range("B1").Value = "vf"
$ws3.Range("I1").Value = "mR"
$ws3.Range("J1").Value = "mB"

$viB3 = @(0.66, 0.618, 0.715, 0.636, 0.398, 0.404, 0.366, 0.309, 0.287, 0.352)
$vf3  = @(0.212, 0.199, 0.231, 0.205, 0.125, 0.128, 0.116, 0.097, 0.09, 0.112)
for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws3.Range("A$r").Value = $viB3[$i]
    $ws3.Range("B$r").Value = $vf3[$i]
}
$ws3.Range("I2").Value = 0.503
$ws3.Range("J2").Value = 0.247

$ws3.Range("A2:B10").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# 4) New worksheet "urto-elastico+massa"
# ---------------------------------------------------------------------------
$lastIdx2 = $wb.Worksheets.Count
$lastSheet2 = $wb.Worksheets.Item($lastIdx2)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$ws4.Name = "urto-elastico+massa"

$ws4.Range("A1").Value = "viB"
$ws4.Range("B1").Value = "vfR"
$ws4.Range("I1").Value = "mR"
$ws4.Range("J1").Value = "mB"
